$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells being written keep their original text
# representation (many look like plain numbers, e.g. "21.99", and Excel
# would otherwise silently coerce them to floating point numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "27.046.55"
$ws.Range("E2").Value2 = "  +0.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.566.29"

$ws.Range("E4").Value2 = "  +0.60%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "208.59"
$ws.Range("E5").Value2 = "  +1.03%  "

$ws.Range("E6").Value2 = "  +0.50%  "

$ws.Range("E7").Value2 = "  +0.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "21.99"
$ws.Range("E8").Value2 = "  +0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.249"
$ws.Range("E9").Value2 = "  +0.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.0597"
$ws.Range("E10").Value2 = "  +0.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.0863"
$ws.Range("E11").Value2 = "  +0.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "1.787.61"
$ws.Range("E12").Value2 = "  +0.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "1.519.30"
$ws.Range("E13").Value2 = "  -1.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "3.78"
$ws.Range("E14").Value2 = "  +0.94%  "

$ws.Range("E15").Value2 = "  +0.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "27.037.23"
$ws.Range("E16").Value2 = "  +0.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "62.03"
$ws.Range("E17").Value2 = "  +0.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.0₃0705"
$ws.Range("E18").Value2 = "  -1.09%  "

$ws.Range("E19").Value2 = "  +1.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "215.13"
$ws.Range("E20").Value2 = "  -0.73%  "

$ws.Range("E21").Value2 = "  +0.69%  "

$ws.Range("E22").Value2 = "  +1.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "9.20"
$ws.Range("E23").Value2 = "  +0.04%  "

$ws.Range("E24").Value2 = "  -0.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "153.99"
$ws.Range("E25").Value2 = "  +0.48%  "

$ws.Range("E26").Value2 = "  -0.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "15.09"
$ws.Range("E27").Value2 = "  +0.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "0.106"
$ws.Range("E28").Value2 = "  +1.29%  "

$ws.Range("E29").Value2 = "  +0.82%  "

$ws.Range("E30").Value2 = "  +0.99%  "

$ws.Range("E31").Value2 = "  +4.46%  "

$ws.Range("E32").Value2 = "  +0.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "3.18"
$ws.Range("E33").Value2 = "  +2.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "1.434.55"
$ws.Range("E34").Value2 = "  +2.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.13"
$ws.Range("E35").Value2 = "  +17.93%  "

$ws.Range("E36").Value2 = "  +1.26%  "

$ws.Range("E37").Value2 = "  +3.06%  "

$ws.Range("E38").Value2 = "  +1.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.532"
$ws.Range("E39").Value2 = "  +1.27%  "

$ws.Range("E40").Value2 = "  +3.07%  "

$ws.Range("E41").Value2 = "  +0.44%  "

$ws.Range("E42").Value2 = "  +4.04%  "

$ws.Range("E43").Value2 = "  +0.82%  "

$ws.Range("E44").Value2 = "  -0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "64.61"
$ws.Range("E45").Value2 = "  +0.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "1.74"
$ws.Range("E46").Value2 = "  +0.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "1.705.93"
$ws.Range("E47").Value2 = "  +1.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "85.84"
$ws.Range("E48").Value2 = "  -1.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.0₆0103"
$ws.Range("E49").Value2 = "  +1.84%  "

$ws.Range("E50").Value2 = "  -0.26%  "

$ws.Range("E51").Value2 = "  +0.54%  "
